# Update recalculated report totals (source data path switched to OneDrive).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 482489.42
$ws.Range("M4").Value = 482489.42
$ws.Range("H5").Value = 482489.42
$ws.Range("N5").Value = 482489.42
$ws.Range("G6").Value = 49090573.21
$ws.Range("M8").Value = 2501127.79
$ws.Range("M10").Value = 46589445.42
$ws.Range("H12").Value = 49090573.21
$ws.Range("N12").Value = 49090573.21
$ws.Range("G13").Value = 3932694.17
$ws.Range("M15").Value = 3932694.17
$ws.Range("H17").Value = 3932694.17
$ws.Range("N17").Value = 3932694.17
$ws.Range("G18").Value = 415427.35
$ws.Range("M20").Value = 107676.85
$ws.Range("M22").Value = 307750.5
$ws.Range("H23").Value = 415427.35
$ws.Range("N23").Value = 415427.35
$ws.Range("G24").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("G27").Value = 3641846.54
$ws.Range("M29").Value = 3641846.54
$ws.Range("H30").Value = 3641846.54
$ws.Range("N30").Value = 3641846.54
$ws.Range("G31").Value = 388882.66
$ws.Range("M33").Value = 388882.66
$ws.Range("H34").Value = 388882.66
$ws.Range("N34").Value = 388882.66
$ws.Range("G35").Value = 195310.33
$ws.Range("M37").Value = 195310.33
$ws.Range("H38").Value = 195310.33
$ws.Range("N38").Value = 195310.33
$ws.Range("G39").Value = 388882.66
$ws.Range("M41").Value = 388882.66
$ws.Range("H42").Value = 388882.66
$ws.Range("N42").Value = 388882.66
$ws.Range("G43").Value = 542214.83
$ws.Range("M45").Value = 542214.83
$ws.Range("H46").Value = 542214.83
$ws.Range("N46").Value = 542214.83
